$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 14 (CB28), shifting CB28..sumPCB down by 2 rows.
$ws.Rows.Item(14).Resize(2).Insert()

# Fill in the new row 14 (CB187) data.
$ws.Range("A14").Value = "CB187"
$ws.Range("B14").Value = 507.98
$ws.Range("C14").Value = 599.96
$ws.Range("D14").Value = 656.14
$ws.Range("E14").Value = 259.29
$ws.Range("F14").Value = 374.96
$ws.Range("G14").Value = 744.5700000000001
$ws.Range("H14").Value = 0.33
$ws.Range("I14").Value = 3.97
$ws.Range("J14").Value = 15.24

# Fill in the new row 15 (CB194) data.
$ws.Range("A15").Value = "CB194"
$ws.Range("B15").Value = 79.56999999999999
$ws.Range("C15").Value = 90.93000000000001
$ws.Range("D15").Value = 104.4
$ws.Range("E15").Value = 53.84
$ws.Range("F15").Value = 84.40000000000001
$ws.Range("G15").Value = 116.23
$ws.Range("H15").Value = 0.33
$ws.Range("I15").Value = 1.24
$ws.Range("J15").Value = 6.98
